$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Replace the paragraph that begins "En este caso sabemos..." with the new
#    text, which includes an inline m:oMath equation (Log2(1177)).
# ---------------------------------------------------------------------------

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.StartsWith("En este caso sabemos")) {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $xml = @"
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math">
        <w:body>
          <w:p w14:paraId="38113B99" w14:textId="268E2D5B" w:rsidR="00D27906" w:rsidRPr="00AA0505" w:rsidRDefault="00D27906" w:rsidP="0063268C" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
            <w:pPr>
              <w:pStyle w:val="Prrafodelista"/>
              <w:spacing w:after="0"/>
              <w:jc w:val="both"/>
              <w:rPr>
                <w:rFonts w:ascii="Dax-Regular" w:hAnsi="Dax-Regular"/>
                <w:lang w:val="es-CO"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Dax-Regular" w:hAnsi="Dax-Regular"/>
                <w:lang w:val="es-CO"/>
              </w:rPr>
              <w:t>E</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Dax-Regular" w:hAnsi="Dax-Regular"/>
                <w:lang w:val="es-CO"/>
              </w:rPr>
              <w:t xml:space="preserve">n este caso podemos ver que este árbol no es lleno, ni completo, ni balanceado. Ya que si calculamos </w:t>
            </w:r>
            <m:oMath><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:lang w:val="es-CO"/></w:rPr><m:t>Log2(1177)</m:t></m:r></m:oMath>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Dax-Regular" w:hAnsi="Dax-Regular"/>
                <w:lang w:val="es-CO"/>
              </w:rPr>
              <w:t xml:space="preserve"> este debería dar la altura del árbol para estar balanceado, y este da 10. Lo que nos lleva a la conclusión que este árbol no se encuentra en balance, ya que su altura es 29. </w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@
    $target.Range.InsertXML($xml)
}

Write-Host "Done."
